$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exercicio 8 (row 9) and Exercicio 9 (row 10) are now done -> mark boolean column C as TRUE
$ws.Range("C9").Value = $true
$ws.Range("C10").Value = $true

# Update the active selection to B11
$null = $ws.Range("B11").Select()
